# "wrap up session summary 1.6.21"
# Adds a new "session summary" note to the תיעוד (documentation) column
# for the three tasks that were worked on in this session, and tidies up
# the sheet view (zoom/selection), default font, and column layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New session-summary note, applied to the three rows touched in this
#     work session (replaces the previous "אין" placeholder in column G). ---
$summary = "סיכום סשן עבודה - 1.6.21"
$ws.Range("G6").Value = $summary
$ws.Range("G7").Value = $summary
$ws.Range("G8").Value = $summary

# --- Sheet view: zoom back out and move the selection to where work left off. ---
$excel.ActiveWindow.Zoom = 115
$ws.Range("G16").Select() | Out-Null

# --- Unhide the helper columns (B/C/D/G) and refresh their widths now that
#     they're visible again. ---
$ws.Columns("B").Hidden = $false
$ws.Columns("B").ColumnWidth = 10

$ws.Columns("C").Hidden = $false
$ws.Columns("C").ColumnWidth = 17.21875

$ws.Columns("D").Hidden = $false
$ws.Columns("D").ColumnWidth = 8.88671875

$ws.Columns("E").ColumnWidth = 41.44140625

$ws.Columns("G").Hidden = $false
$ws.Columns("G").ColumnWidth = 30.88671875

# --- Switch the workbook's default/Normal font from Arial to Calibri. ---
$wb.Styles.Item("Normal").Font.Name = "Calibri"
